$d = $word.ActiveDocument

# --- Text replacements (run text content) ---
$found0 = $d.Content.Find.Execute("This document contains all homework assignments for Group 2 from week 8-15.", $true, $false, $false, $false, $false, $true, 1, $false, "This document contains all homework assignments for Group 2 from week 8-15. Participants of this group include Vinicio Haro, Juliann McEachern, Jeremy O’Brien, Bethany Poulin, and Sang (Andy) Yoon.", 2)
if (-not $found0) { Write-Output "WARNING: text replacement 0 not found" }
$found1 = $d.Content.Find.Execute("(b) A small percentage of cells in the predictor set contain missing values. Use an imputation function to ï¬ll in these missing values (e.g., see Sect. 3.8).", $true, $false, $false, $false, $false, $true, 1, $false, "(b) A small percentage of cells in the predictor set contain missing values. Use an imputation function to fill in these missing values (e.g., see Sect. 3.8).", 2)
if (-not $found1) { Write-Output "WARNING: text replacement 1 not found" }
$found2 = $d.Content.Find.Execute("(b). Which models appear to give the best performance? Does MARS select the informative predictors (those named X1âX5)?", $true, $false, $false, $false, $false, $true, 1, $false, "(b). Which models appear to give the best performance? Does MARS select the informative predictors (those named X1-X5)?", 2)
if (-not $found2) { Write-Output "WARNING: text replacement 2 not found" }
$found3 = $d.Content.Find.Execute("(a). Fit a random forest model to all of the predictors, then estimate the variable importance scores. Did the random forest model signiï¬cantly use the uninformative predictors (V6 â V10)?", $true, $false, $false, $false, $false, $true, 1, $false, "(a). Fit a random forest model to all of the predictors, then estimate the variable importance scores. Did the random forest model significantly use the uninformative predictors (V6-V10)?", 2)
if (-not $found3) { Write-Output "WARNING: text replacement 3 not found" }
$found4 = $d.Content.Find.Execute("function in the party package to ï¬t a random forest model using conditional inference trees. The party package function", $true, $false, $false, $false, $false, $true, 1, $false, "function in the party package to fit a random forest model using conditional inference trees. The party package function", 2)
if (-not $found4) { Write-Output "WARNING: text replacement 4 not found" }
$found5 = $d.Content.Find.Execute("argument of that function toggles between the traditional importance measure and the modiï¬ed version described in Strobl et al. (2007). Do these importances show the same pattern as the traditional random forest model?", $true, $false, $false, $false, $false, $true, 1, $false, "argument of that function toggles between the traditional importance measure and the modified version described in Strobl et al. (2007). Do these importances show the same pattern as the traditional random forest model?", 2)
if (-not $found5) { Write-Output "WARNING: text replacement 5 not found" }
$found6 = $d.Content.Find.Execute("(d). Repeat this process with diï¬erent tree models, such as boosted trees and Cubist. Does the same pattern occur?", $true, $false, $false, $false, $false, $true, 1, $false, "(d). Repeat this process with different tree models, such as boosted trees and Cubist. Does the same pattern occur?", 2)
if (-not $found6) { Write-Output "WARNING: text replacement 6 not found" }
$found7 = $d.Content.Find.Execute("8.2: Use a simulation to show tree bias with diï¬erent granularities.", $true, $false, $false, $false, $false, $true, 1, $false, "8.2: Use a simulation to show tree bias with different granularities.", 2)
if (-not $found7) { Write-Output "WARNING: text replacement 7 not found" }
$found8 = $d.Content.Find.Execute("8.3: In stochastic gradient boosting the bagging fraction and learning rate will govern the construction of the trees as they are guided by the gradient. Although the optimal values of these parameters should be obtained through the tuning process, it is helpful to understand how the magnitudes of these parameters aï¬ect magnitudes of variable importance. Figure 8.24 provides the variable importance plots for boosting using two extreme values for the bagging fraction (0.1 and 0.9) and the learning rate (0.1 and 0.9) for the solubility data. The left-hand plot has both parameters set to 0.1, and the right-hand plot has both set to 0.9:", $true, $false, $false, $false, $false, $true, 1, $false, "8.3: In stochastic gradient boosting the bagging fraction and learning rate will govern the construction of the trees as they are guided by the gradient. Although the optimal values of these parameters should be obtained through the tuning process, it is helpful to understand how the magnitudes of these parameters affect magnitudes of variable importance. Figure 8.24 provides the variable importance plots for boosting using two extreme values for the bagging fraction (0.1 and 0.9) and the learning rate (0.1 and 0.9) for the solubility data. The left-hand plot has both parameters set to 0.1, and the right-hand plot has both set to 0.9:", 2)
if (-not $found8) { Write-Output "WARNING: text replacement 8 not found" }
$found9 = $d.Content.Find.Execute("(a). Why does the model on the right focus its importance on just the ï¬rst few of predictors, whereas the model on the left spreads importance across more predictors?", $true, $false, $false, $false, $false, $true, 1, $false, "(a). Why does the model on the right focus its importance on just the first few of predictors, whereas the model on the left spreads importance across more predictors?", 2)
if (-not $found9) { Write-Output "WARNING: text replacement 9 not found" }
$found10 = $d.Content.Find.Execute("(c). How would increasing interaction depth aï¬ect the slope of predictor importance for either model in Fig.8.24?", $true, $false, $false, $false, $false, $true, 1, $false, "(c). How would increasing interaction depth affect the slope of predictor importance for either model in Fig.8.24?", 2)
if (-not $found10) { Write-Output "WARNING: text replacement 10 not found" }

# --- Bookmark renames (preserve range, rename bookmark) ---
if ($d.Bookmarks.Exists("b-a-small-percentage-of-cells-in-the-predictor-set-contain-missing-values.-use-an-imputation-function-to-ill-in-these-missing-values-e.g.-see-sect.-3.8.")) {
    $bm0 = $d.Bookmarks.Item("b-a-small-percentage-of-cells-in-the-predictor-set-contain-missing-values.-use-an-imputation-function-to-ill-in-these-missing-values-e.g.-see-sect.-3.8.")
    $rng0 = $bm0.Range
    $bm0.Delete()
    $d.Bookmarks.Add("b-a-small-percentage-of-cells-in-the-predictor-set-contain-missing-values.-use-an-imputation-function-to-fill-in-these-missing-values-e.g.-see-sect.-3.8.", $rng0) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: b-a-small-percentage-of-cells-in-the-predictor-set-contain-missing-values.-use-an-imputation-function-to-ill-in-these-missing-values-e.g.-see-sect.-3.8."
}
if ($d.Bookmarks.Exists("b.-which-models-appear-to-give-the-best-performance-does-mars-select-the-informative-predictors-those-named-x1ax5")) {
    $bm1 = $d.Bookmarks.Item("b.-which-models-appear-to-give-the-best-performance-does-mars-select-the-informative-predictors-those-named-x1ax5")
    $rng1 = $bm1.Range
    $bm1.Delete()
    $d.Bookmarks.Add("b.-which-models-appear-to-give-the-best-performance-does-mars-select-the-informative-predictors-those-named-x1-x5", $rng1) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: b.-which-models-appear-to-give-the-best-performance-does-mars-select-the-informative-predictors-those-named-x1ax5"
}
if ($d.Bookmarks.Exists("a.-fit-a-random-forest-model-to-all-of-the-predictors-then-estimate-the-variable-importance-scores.-did-the-random-forest-model-signiicantly-use-the-uninformative-predictors-v6-a-v10")) {
    $bm2 = $d.Bookmarks.Item("a.-fit-a-random-forest-model-to-all-of-the-predictors-then-estimate-the-variable-importance-scores.-did-the-random-forest-model-signiicantly-use-the-uninformative-predictors-v6-a-v10")
    $rng2 = $bm2.Range
    $bm2.Delete()
    $d.Bookmarks.Add("a.-fit-a-random-forest-model-to-all-of-the-predictors-then-estimate-the-variable-importance-scores.-did-the-random-forest-model-significantly-use-the-uninformative-predictors-v6-v10", $rng2) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: a.-fit-a-random-forest-model-to-all-of-the-predictors-then-estimate-the-variable-importance-scores.-did-the-random-forest-model-signiicantly-use-the-uninformative-predictors-v6-a-v10"
}
if ($d.Bookmarks.Exists("c.-use-the-cforest-function-in-the-party-package-to-it-a-random-forest-model-using-conditional-inference-trees.-the-party-package-function-varimp-can-calculate-predictor-importance.-the-conditional-argument-of-that-function-toggles-between-the-traditional-importance-measure-and-the-modiied-version-described-in-strobl-et-al.-2007.-do-these-importances-show-the-same-pattern-as-the-traditional-random-forest-model")) {
    $bm3 = $d.Bookmarks.Item("c.-use-the-cforest-function-in-the-party-package-to-it-a-random-forest-model-using-conditional-inference-trees.-the-party-package-function-varimp-can-calculate-predictor-importance.-the-conditional-argument-of-that-function-toggles-between-the-traditional-importance-measure-and-the-modiied-version-described-in-strobl-et-al.-2007.-do-these-importances-show-the-same-pattern-as-the-traditional-random-forest-model")
    $rng3 = $bm3.Range
    $bm3.Delete()
    $d.Bookmarks.Add("c.-use-the-cforest-function-in-the-party-package-to-fit-a-random-forest-model-using-conditional-inference-trees.-the-party-package-function-varimp-can-calculate-predictor-importance.-the-conditional-argument-of-that-function-toggles-between-the-traditional-importance-measure-and-the-modified-version-described-in-strobl-et-al.-2007.-do-these-importances-show-the-same-pattern-as-the-traditional-random-forest-model", $rng3) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: c.-use-the-cforest-function-in-the-party-package-to-it-a-random-forest-model-using-conditional-inference-trees.-the-party-package-function-varimp-can-calculate-predictor-importance.-the-conditional-argument-of-that-function-toggles-between-the-traditional-importance-measure-and-the-modiied-version-described-in-strobl-et-al.-2007.-do-these-importances-show-the-same-pattern-as-the-traditional-random-forest-model"
}
if ($d.Bookmarks.Exists("d.-repeat-this-process-with-diierent-tree-models-such-as-boosted-trees-and-cubist.-does-the-same-pattern-occur")) {
    $bm4 = $d.Bookmarks.Item("d.-repeat-this-process-with-diierent-tree-models-such-as-boosted-trees-and-cubist.-does-the-same-pattern-occur")
    $rng4 = $bm4.Range
    $bm4.Delete()
    $d.Bookmarks.Add("d.-repeat-this-process-with-different-tree-models-such-as-boosted-trees-and-cubist.-does-the-same-pattern-occur", $rng4) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: d.-repeat-this-process-with-diierent-tree-models-such-as-boosted-trees-and-cubist.-does-the-same-pattern-occur"
}
if ($d.Bookmarks.Exists("use-a-simulation-to-show-tree-bias-with-diierent-granularities.")) {
    $bm5 = $d.Bookmarks.Item("use-a-simulation-to-show-tree-bias-with-diierent-granularities.")
    $rng5 = $bm5.Range
    $bm5.Delete()
    $d.Bookmarks.Add("use-a-simulation-to-show-tree-bias-with-different-granularities.", $rng5) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: use-a-simulation-to-show-tree-bias-with-diierent-granularities."
}
if ($d.Bookmarks.Exists("in-stochastic-gradient-boosting-the-bagging-fraction-and-learning-rate-will-govern-the-construction-of-the-trees-as-they-are-guided-by-the-gradient.-although-the-optimal-values-of-these-parameters-should-be-obtained-through-the-tuning-process-it-is-helpful-to-understand-how-the-magnitudes-of-these-parameters-aiect-magnitudes-of-variable-importance.-figure-8.24-provides-the-variable-importance-plots-for-boosting-using-two-extreme-values-for-the-bagging-fraction-0.1-and-0.9-and-the-learning-rate-0.1-and-0.9-for-the-solubility-data.-the-left-hand-plot-has-both-parameters-set-to-0.1-and-the-right-hand-plot-has-both-set-to-0.9")) {
    $bm6 = $d.Bookmarks.Item("in-stochastic-gradient-boosting-the-bagging-fraction-and-learning-rate-will-govern-the-construction-of-the-trees-as-they-are-guided-by-the-gradient.-although-the-optimal-values-of-these-parameters-should-be-obtained-through-the-tuning-process-it-is-helpful-to-understand-how-the-magnitudes-of-these-parameters-aiect-magnitudes-of-variable-importance.-figure-8.24-provides-the-variable-importance-plots-for-boosting-using-two-extreme-values-for-the-bagging-fraction-0.1-and-0.9-and-the-learning-rate-0.1-and-0.9-for-the-solubility-data.-the-left-hand-plot-has-both-parameters-set-to-0.1-and-the-right-hand-plot-has-both-set-to-0.9")
    $rng6 = $bm6.Range
    $bm6.Delete()
    $d.Bookmarks.Add("in-stochastic-gradient-boosting-the-bagging-fraction-and-learning-rate-will-govern-the-construction-of-the-trees-as-they-are-guided-by-the-gradient.-although-the-optimal-values-of-these-parameters-should-be-obtained-through-the-tuning-process-it-is-helpful-to-understand-how-the-magnitudes-of-these-parameters-affect-magnitudes-of-variable-importance.-figure-8.24-provides-the-variable-importance-plots-for-boosting-using-two-extreme-values-for-the-bagging-fraction-0.1-and-0.9-and-the-learning-rate-0.1-and-0.9-for-the-solubility-data.-the-left-hand-plot-has-both-parameters-set-to-0.1-and-the-right-hand-plot-has-both-set-to-0.9", $rng6) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: in-stochastic-gradient-boosting-the-bagging-fraction-and-learning-rate-will-govern-the-construction-of-the-trees-as-they-are-guided-by-the-gradient.-although-the-optimal-values-of-these-parameters-should-be-obtained-through-the-tuning-process-it-is-helpful-to-understand-how-the-magnitudes-of-these-parameters-aiect-magnitudes-of-variable-importance.-figure-8.24-provides-the-variable-importance-plots-for-boosting-using-two-extreme-values-for-the-bagging-fraction-0.1-and-0.9-and-the-learning-rate-0.1-and-0.9-for-the-solubility-data.-the-left-hand-plot-has-both-parameters-set-to-0.1-and-the-right-hand-plot-has-both-set-to-0.9"
}
if ($d.Bookmarks.Exists("a.-why-does-the-model-on-the-right-focus-its-importance-on-just-the-irst-few-of-predictors-whereas-the-model-on-the-left-spreads-importance-across-more-predictors")) {
    $bm7 = $d.Bookmarks.Item("a.-why-does-the-model-on-the-right-focus-its-importance-on-just-the-irst-few-of-predictors-whereas-the-model-on-the-left-spreads-importance-across-more-predictors")
    $rng7 = $bm7.Range
    $bm7.Delete()
    $d.Bookmarks.Add("a.-why-does-the-model-on-the-right-focus-its-importance-on-just-the-first-few-of-predictors-whereas-the-model-on-the-left-spreads-importance-across-more-predictors", $rng7) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: a.-why-does-the-model-on-the-right-focus-its-importance-on-just-the-irst-few-of-predictors-whereas-the-model-on-the-left-spreads-importance-across-more-predictors"
}
if ($d.Bookmarks.Exists("c.-how-would-increasing-interaction-depth-aiect-the-slope-of-predictor-importance-for-either-model-in-fig.8.24")) {
    $bm8 = $d.Bookmarks.Item("c.-how-would-increasing-interaction-depth-aiect-the-slope-of-predictor-importance-for-either-model-in-fig.8.24")
    $rng8 = $bm8.Range
    $bm8.Delete()
    $d.Bookmarks.Add("c.-how-would-increasing-interaction-depth-affect-the-slope-of-predictor-importance-for-either-model-in-fig.8.24", $rng8) | Out-Null
} else {
    Write-Output "WARNING: bookmark not found: c.-how-would-increasing-interaction-depth-aiect-the-slope-of-predictor-importance-for-either-model-in-fig.8.24"
}

Write-Output "done"
